$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.006.47"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.256.64"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'306.38"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'96.53"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").Value = "'0.523"
$ws.Range("E7").Value = "  -1.46%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "'34.82"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "'6.80"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "2.606.16"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "2.268.25"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'0.785"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").Value = "41.878.65"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'12.17"
$ws.Range("E19").Value = "  -4.01%  "
$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "'67.46"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").Value = "'235.70"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").Value = "'36.54"
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.12"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.51"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").Value = "'164.46"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("D35").Value = "'17.51"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("E38").Value = "  -5.11%  "
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").Value = "'1.80"
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("D41").Value = "'4.10"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("D42").Value = "'2.33"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").Value = "1.947.59"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").Value = "'18.93"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("D46").Value = "'10.00"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").Value = "'53.06"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "2.479.80"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'91.88"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'71.56"
$ws.Range("E51").Value = "  -1.50%  "
